$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item(1)
$ws.Range("F2").Value = 232
$ws.Range("F3").Value = 375
$ws.Range("F5").Value = 24
$ws.Range("F6").Value = 91
$ws.Range("F9").Value = 4635
$ws.Range("F10").Value = 4635
$ws.Range("F12").Value = 446
$ws.Range("F13").Value = 1084
$ws.Range("F14").Value = 602
$ws.Range("F15").Value = 4131
$ws.Range("F16").Value = 159
$ws.Range("F17").Value = 160
$ws.Range("F18").Value = 47
$ws.Range("F19").Value = 202
$ws.Range("F20").Value = 3415
$ws.Range("F24").Value = 2965
$ws.Range("F25").Value = 122
$ws.Range("F26").Value = 122
$ws.Range("F28").Value = 142
$ws.Range("F29").Value = 178
$ws.Range("F30").Value = 172
$ws.Range("F31").Value = 71
$ws.Range("F32").Value = 48
$ws.Range("F36").Value = 5356
$ws.Range("F37").Value = 739
$ws.Range("F38").Value = 383
$ws.Range("F39").Value = 80
$ws.Range("F41").Value = 27
$ws.Range("F42").Value = 1059
$ws.Range("F43").Value = 439
$ws.Range("F45").Value = 1940
$ws.Range("F46").Value = 295
$ws.Range("F48").Value = 688
$ws.Range("F49").Value = 826

$ws = $wb.Worksheets.Item(2)
$ws.Range("F6").Value = 79
$ws.Range("F22").Value = 719

$ws = $wb.Worksheets.Item(3)
$ws.Range("F2").Value = 196

$ws = $wb.Worksheets.Item(4)
$ws.Range("F2").Value = 196
$ws.Range("F5").Value = 232
$ws.Range("F6").Value = 24
$ws.Range("F7").Value = 79
$ws.Range("F8").Value = 91
$ws.Range("F11").Value = 4635
$ws.Range("F12").Value = 4635
$ws.Range("F17").Value = 446
$ws.Range("F18").Value = 1084
$ws.Range("F19").Value = 602
$ws.Range("F20").Value = 4131
$ws.Range("F21").Value = 159
$ws.Range("F22").Value = 160
$ws.Range("F23").Value = 202
$ws.Range("F24").Value = 3415
$ws.Range("F25").Value = 2965
$ws.Range("F26").Value = 122
$ws.Range("F27").Value = 122
$ws.Range("F28").Value = 142
$ws.Range("F29").Value = 178
$ws.Range("F30").Value = 172
$ws.Range("F37").Value = 5356
$ws.Range("F39").Value = 739
$ws.Range("F40").Value = 383
$ws.Range("F42").Value = 80
$ws.Range("F44").Value = 1059
$ws.Range("F45").Value = 439
$ws.Range("F47").Value = 1940
$ws.Range("F49").Value = 688
$ws.Range("F50").Value = 826

Write-Output "done"